# redmine 4180 - documenting 'attibute' fields from eml
#
# Adds 5 new rows (SOLR field definitions) describing the EML
# <attribute> related fields to the "EML" worksheet, then makes "EML"
# the active/selected sheet (it was "Dryad" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EML")

# --- New field rows (A:D) -------------------------------------------------
$ws.Cells.Item(50, 1).Value = "attributeName"
$ws.Cells.Item(50, 2).Value = "string"
$ws.Cells.Item(50, 3).Value = "Yes"
$ws.Cells.Item(50, 4).Value = "//dataTable/attributeList/attribute/attributeName/text()"

$ws.Cells.Item(51, 1).Value = "attributeLabel"
$ws.Cells.Item(51, 2).Value = "string"
$ws.Cells.Item(51, 3).Value = "Yes"
$ws.Cells.Item(51, 4).Value = "//dataTable/attributeList/attribute/attributeLabel/text()"

$ws.Cells.Item(52, 1).Value = "attributeDescription"
$ws.Cells.Item(52, 2).Value = "string"
$ws.Cells.Item(52, 3).Value = "Yes"
$ws.Cells.Item(52, 4).Value = "//dataTable/attributeList/attribute/attributeDefinition/text()"

$ws.Cells.Item(53, 1).Value = "attributeUnit"
$ws.Cells.Item(53, 2).Value = "string"
$ws.Cells.Item(53, 3).Value = "Yes"
$ws.Cells.Item(53, 4).Value = "//dataTable//standardUnit/text() | //dataTable//customUnit/text()"

$ws.Cells.Item(54, 1).Value = "attribute"
$ws.Cells.Item(54, 2).Value = "string "
$ws.Cells.Item(54, 3).Value = "Yes"
$ws.Cells.Item(54, 4).Value = "//dataTable/attributeList/attribute"

# --- View state: EML becomes the active tab/selection ---------------------
$ws.Activate()
$ws.Range("A55").Select()
